$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the "_GoBack" bookmark currently sits right after "... SNGF" (end
# of the "Ha mai fatto un ritaglio?" bullet). Move it so it instead sits
# inside the "DOMANDE INIZIALI" heading, right after "DOMANDE INIZ" (this
# splits that heading's single run into "DOMANDE INIZ" + bookmark + "IALI",
# matching the target OOXML).
# ---------------------------------------------------------------------------

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$headingHit = $d.Content
$headingHit.Find.Execute("DOMANDE INIZ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$splitPoint = $headingHit.End
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint))

# ---------------------------------------------------------------------------
# Change 2: in the "Cosa si intende fare..." bullet, the trailing "?" and
# " " are currently two separate runs; collapse them into a single run
# whose text is "? " (the preceding "Cosa si intende..." run is untouched).
# A plain Range.Text assignment would get silently re-merged into the
# preceding run by the engine's run-coalescing pass (same, empty,
# formatting on both sides), so the replacement is done with InsertXML,
# scoped exactly to the "?" + " " range, which keeps it as its own run.
# ---------------------------------------------------------------------------

$qHit = $d.Content
$qHit.Find.Execute("Cosa si intende fare: pavimento, rivestimento o entrambi? ", `
                    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$questionMarkSpace = $d.Range($qHit.End - 2, $qHit.End)

$questionMarkSpace.InsertXML(@"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml"
            pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">? </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@)
